$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the saved window view (best-effort; mirrors the author's
# recorded workbookView window dimensions in the target file).
try {
    $excel.ActiveWindow.Width = 20490
    $excel.ActiveWindow.Height = 7545
} catch {}

$ws.Range("D1").Value = "ORG_TREA_IDENOLD"
$ws.Range("E1").Value = "ORG_TREA_IDENNEW"
$ws.Range("F1").Value = "ORG_TREA_STATUS"

$ws.Range("D1").HorizontalAlignment = -4131
$ws.Range("F1").HorizontalAlignment = -4131

$ws.Range("F7").Select()
